$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 338, shifting rows 338-362
# down to 340-364 (their original data is preserved by Excel's Insert).
$ws.Rows("338:339").Insert()

# Populate the two newly inserted rows with the new record data.
$ws.Range("A338").Value = 10
$ws.Range("B338").Value = "Vega Modelo de Temuco"
$ws.Range("C338").Value = "La Araucanía"
$ws.Range("D338").Value = 44461
$ws.Range("E338").Value = 9
$ws.Range("F338").Value = 100112003
$ws.Range("G338").Value = "Ajo"
$ws.Range("H338").Value = "Chino"
$ws.Range("I338").Value = "Primera"
$ws.Range("J338").Value = 100
$ws.Range("K338").Value = 18000
$ws.Range("L338").Value = 18000
$ws.Range("M338").Value = 18000
$ws.Range("N338").Value = "$/caja 10 kilos"
$ws.Range("O338").Value = "China"
$ws.Range("P338").Value = 1800
$ws.Range("Q338").Value = 10
$ws.Range("R338").Value = "Hortaliza"

$ws.Range("A339").Value = 10
$ws.Range("B339").Value = "Vega Modelo de Temuco"
$ws.Range("C339").Value = "La Araucanía"
$ws.Range("D339").Value = 44461
$ws.Range("E339").Value = 9
$ws.Range("F339").Value = 100112003
$ws.Range("G339").Value = "Ajo"
$ws.Range("H339").Value = "Chino"
$ws.Range("I339").Value = "Primera"
$ws.Range("J339").Value = 220
$ws.Range("K339").Value = 19000
$ws.Range("L339").Value = 20000
$ws.Range("M339").Value = 19545
$ws.Range("N339").Value = "$/malla 10 kilos"
$ws.Range("O339").Value = "China"
$ws.Range("P339").Value = 1954
$ws.Range("Q339").Value = 10
$ws.Range("R339").Value = "Hortaliza"
